$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

# The source data in this sheet stores every value (even numeric-looking
# ones like "9.4" or "88673") as plain text, so force Text format on the
# cells we touch to keep them stored as text rather than being coerced to
# numbers on write.
$ws.Range("A11:D12").NumberFormat = "@"
$ws.Range("A33:D34").NumberFormat = "@"

# Table 1 (rows 11-12): swap "Enterprises (absolute #)" row with
# "Enterprises density (per 1000 people)" row so the density row comes first.
$a11 = [string]$ws.Range("A11").Value()
$d11 = [string]$ws.Range("D11").Value()
$a12 = [string]$ws.Range("A12").Value()
$d12 = [string]$ws.Range("D12").Value()

$ws.Range("A11").Value = $a12
$ws.Range("D11").Value = $d12
$ws.Range("A12").Value = $a11
$ws.Range("D12").Value = $d11

# Table 3 (rows 33-34): same swap, this time across columns A, B, C, D.
$a33 = [string]$ws.Range("A33").Value()
$b33 = [string]$ws.Range("B33").Value()
$c33 = [string]$ws.Range("C33").Value()
$d33 = [string]$ws.Range("D33").Value()
$a34 = [string]$ws.Range("A34").Value()
$b34 = [string]$ws.Range("B34").Value()
$c34 = [string]$ws.Range("C34").Value()
$d34 = [string]$ws.Range("D34").Value()

$ws.Range("A33").Value = $a34
$ws.Range("B33").Value = $b34
$ws.Range("C33").Value = $c34
$ws.Range("D33").Value = $d34
$ws.Range("A34").Value = $a33
$ws.Range("B34").Value = $b33
$ws.Range("C34").Value = $c33
$ws.Range("D34").Value = $d33
